# Rename sheet 'Data' to 'Data table' (Close #151)

$wb = $excel.ActiveWorkbook

# Last selection on the "Table" sheet (captured before focus moves away).
$tableSheet = $wb.Worksheets.Item("Table")
[void]$tableSheet.Range("B89").Select()

# Rename "Data" -> "Data table".
$dataSheet = $wb.Worksheets.Item("Data")
$dataSheet.Name = "Data table"

# "Data table" becomes the active sheet.
$dataSheet.Activate()
